$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Donor cell with the standard data-cell style (border-only, General format)
# used to restore formatting on cells where assigning a "NN%" text value
# forces Excel to re-type the cell (percentage number + new style).
$styleDonor = $ws.Range("H2")

$ws.Range("E2").Value = "2026-02-10 05:48:17"
$ws.Range("E3").Value = "2026-02-10 05:48:19"
$ws.Range("G3").Value = "182 cm"
$ws.Range("I3").Value = "7.3 mm"
$ws.Range("E4").Value = "2026-02-10 05:48:22"
$styleDonor.Copy()
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "87%"
$ws.Range("H4").PasteSpecial(-4122)
$ws.Range("N4").Value = "7.0 °C 5:09 TU"
$ws.Range("O4").Value = "8.7 °C"
$ws.Range("E5").Value = "2026-02-10 05:48:24"
$ws.Range("G5").Value = "131 cm"
$ws.Range("I5").Value = "11.6 mm"
$ws.Range("E6").Value = "2026-02-10 05:48:26"
$ws.Range("N6").Value = "6.3 °C 5:18 TU"
$ws.Range("E7").Value = "2026-02-10 05:48:29"
$styleDonor.Copy()
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "79%"
$ws.Range("H7").PasteSpecial(-4122)
$ws.Range("L7").Value = "20.9 km/h - 304º 5:05 TU"
$ws.Range("E8").Value = "2026-02-10 05:48:31"
$ws.Range("E9").Value = "2026-02-10 05:48:33"
$ws.Range("O9").Value = "6.6 °C"
$ws.Range("E10").Value = "2026-02-10 05:48:36"
$ws.Range("N10").Value = "5.5 °C 5:29 TU"
$ws.Range("O10").Value = "7.2 °C"
$ws.Range("E11").Value = "2026-02-10 05:48:38"
$ws.Range("E12").Value = "2026-02-10 05:48:40"
$ws.Range("N12").Value = "5.3 °C 5:27 TU"
$ws.Range("O12").Value = "6.9 °C"
$ws.Range("E13").Value = "2026-02-10 05:48:43"
$ws.Range("I13").Value = "2.0 mm"
$ws.Range("L13").Value = "9.4 km/h - 67º 5:26 TU"
$ws.Range("N13").Value = "2.4 °C 5:08 TU"
$ws.Range("E14").Value = "2026-02-10 05:48:45"
$ws.Range("O14").Value = "9.3 °C"
$ws.Range("E15").Value = "2026-02-10 05:48:47"
$ws.Range("O15").Value = "6.4 °C"
$ws.Range("E16").Value = "2026-02-10 05:48:50"
$ws.Range("I16").Value = "11.1 mm"
$ws.Range("E17").Value = "2026-02-10 05:48:52"
$styleDonor.Copy()
$ws.Range("H17").NumberFormat = "@"
$ws.Range("H17").Value = "92%"
$ws.Range("H17").PasteSpecial(-4122)
$ws.Range("M17").Value = "5.8 °C 5:14 TU"
$ws.Range("O17").Value = "2.6 °C"
$ws.Range("E18").Value = "2026-02-10 05:48:55"
$ws.Range("N18").Value = "5.3 °C 5:29 TU"
$ws.Range("O18").Value = "7.5 °C"
$ws.Range("E19").Value = "2026-02-10 05:48:57"
$ws.Range("O19").Value = "3.9 °C"
$ws.Range("E20").Value = "2026-02-10 05:48:59"
$ws.Range("E21").Value = "2026-02-10 05:49:02"
$ws.Range("I21").Value = "2.9 mm"
$ws.Range("J21").Value = "1007.5 hPa"
$ws.Range("E22").Value = "2026-02-10 05:49:04"
$ws.Range("G22").Value = "124 cm"
$ws.Range("I22").Value = "0.2 mm"
$ws.Range("E23").Value = "2026-02-10 05:49:06"
$ws.Range("G23").Value = "182 cm"
$ws.Range("I23").Value = "9.5 mm"
$ws.Range("E24").Value = "2026-02-10 05:49:09"
$ws.Range("I24").Value = "1.0 mm"
$ws.Range("E25").Value = "2026-02-10 05:49:11"
$ws.Range("G25").Value = "117 cm"
$ws.Range("I25").Value = "5.9 mm"
$ws.Range("O25").Value = "-0.7 °C"
$ws.Range("E26").Value = "2026-02-10 05:49:14"
$ws.Range("M26").Value = "4.1 °C 5:26 TU"
$ws.Range("E27").Value = "2026-02-10 05:49:16"
$ws.Range("I27").Value = "1.7 mm"
$ws.Range("O27").Value = "-0.4 °C"
$ws.Range("E28").Value = "2026-02-10 05:49:19"
$ws.Range("N28").Value = "3.6 °C 5:13 TU"
$ws.Range("O28").Value = "5.1 °C"
$ws.Range("E29").Value = "2026-02-10 05:49:21"
$ws.Range("N29").Value = "6.6 °C 5:29 TU"
$ws.Range("O29").Value = "8.9 °C"
$ws.Range("E30").Value = "2026-02-10 05:49:23"
$ws.Range("N30").Value = "6.8 °C 5:08 TU"
$ws.Range("O30").Value = "7.4 °C"
$ws.Range("E31").Value = "2026-02-10 05:49:26"
$styleDonor.Copy()
$ws.Range("H31").NumberFormat = "@"
$ws.Range("H31").Value = "85%"
$ws.Range("H31").PasteSpecial(-4122)
$ws.Range("K31").Value = "-0.1 MJ/m2"
$ws.Range("O31").Value = "8.8 °C"
$ws.Range("E32").Value = "2026-02-10 05:49:28"
$ws.Range("E33").Value = "2026-02-10 05:49:31"
$ws.Range("I33").Value = "3.9 mm"
$ws.Range("J33").Value = "1007.8 hPa"
$ws.Range("N33").Value = "1.7 °C 5:21 TU"
$ws.Range("O33").Value = "1.9 °C"
$ws.Range("E34").Value = "2026-02-10 05:49:33"
$ws.Range("I34").Value = "1.9 mm"
$ws.Range("E35").Value = "2026-02-10 05:49:36"
$styleDonor.Copy()
$ws.Range("H35").NumberFormat = "@"
$ws.Range("H35").Value = "82%"
$ws.Range("H35").PasteSpecial(-4122)
$ws.Range("I35").Value = "0.3 mm"
$ws.Range("N35").Value = "9.5 °C 5:09 TU"
$ws.Range("O35").Value = "10.4 °C"
$ws.Range("E36").Value = "2026-02-10 05:49:38"
$styleDonor.Copy()
$ws.Range("H36").NumberFormat = "@"
$ws.Range("H36").Value = "93%"
$ws.Range("H36").PasteSpecial(-4122)
$ws.Range("O36").Value = "9.0 °C"
$ws.Range("E37").Value = "2026-02-10 05:49:40"
$styleDonor.Copy()
$ws.Range("H37").NumberFormat = "@"
$ws.Range("H37").Value = "96%"
$ws.Range("H37").PasteSpecial(-4122)
$ws.Range("J37").Value = "1007.1 hPa"
$ws.Range("O37").Value = "3.8 °C"
$ws.Range("E38").Value = "2026-02-10 05:49:43"
$ws.Range("N38").Value = "6.7 °C 5:20 TU"
$ws.Range("O38").Value = "7.7 °C"
$ws.Range("E39").Value = "2026-02-10 05:49:45"
$ws.Range("I39").Value = "2.3 mm"
$ws.Range("E40").Value = "2026-02-10 05:49:48"
$ws.Range("I40").Value = "3.6 mm"
$ws.Range("J40").Value = "1008.3 hPa"
$ws.Range("N40").Value = "4.1 °C 5:28 TU"
$ws.Range("O40").Value = "4.7 °C"
$ws.Range("E41").Value = "2026-02-10 05:49:50"
$ws.Range("O41").Value = "9.8 °C"
$ws.Range("E42").Value = "2026-02-10 05:49:52"
$ws.Range("O42").Value = "8.1 °C"
$ws.Range("E43").Value = "2026-02-10 05:49:55"
$ws.Range("N43").Value = "5.4 °C 5:17 TU"
$ws.Range("O43").Value = "6.1 °C"
$ws.Range("E44").Value = "2026-02-10 05:49:57"
$ws.Range("I44").Value = "7.1 mm"
$ws.Range("E45").Value = "2026-02-10 05:49:59"
$ws.Range("I45").Value = "14.8 mm"
$ws.Range("O45").Value = "3.2 °C"
$ws.Range("E46").Value = "2026-02-10 05:50:02"
$styleDonor.Copy()
$ws.Range("H46").NumberFormat = "@"
$ws.Range("H46").Value = "99%"
$ws.Range("H46").PasteSpecial(-4122)
$ws.Range("I46").Value = "0.1 mm"
$ws.Range("L46").Value = "10.4 km/h - 52º 5:29 TU"
